# Edit script for LOM3057.docx
# Applies the set of textual changes described by the target diff.

$d = $word.ActiveDocument

function Set-ParagraphRuns {
    # NOTE: named parameter binding is unreliable in this runtime, so this
    # helper relies on positional arguments only.
    param($MatchText, $InnerXml)
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        if ($t -like ("*" + $MatchText + "*")) {
            $full = $p.Range
            # Exclude the trailing paragraph mark from the replaced range.
            $sub = $d.Range($full.Start, $full.End - 1)
            $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                   '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                   '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body>' + $InnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $sub.InsertXML($xml)
            return $true
        }
    }
    return $false
}

# 1) Ativação date: 2020 -> 2024
$d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2024", 2) | Out-Null

# 2) Objetivos (PT) paragraph: replace with the new, longer description.
$oldObjPt = "Introduzir conceitos básicos sobre a influência da estrutura química dos materiais poliméricos na morfologia e nas propriedades mecânicas, ópticas, elétricas e químicas."
$newObjPt = "Abordar conceitos fundamentais sobre materiais poliméricos, envolvendo o histórico de desenvolvimento, organização do setor produtivo, formas de obtenção, estrutura química e física destes materiais e respectiva relação com propriedades térmicas, mecânicas, reológicas, ópticas, elétricas, etc. Permitir que o aluno tenha uma visão clara sobre estrutura, propriedades e aplicações de polímeros termoplásticos, termorrígidos e elastômeros, bem como as propriedades destes materiais podem ser modificadas com o uso de aditivos.Capacitar o aluno com conhecimentos para que possa cursar outras disciplinas na área de materiais poliméricos."
$d.Content.Find.Execute($oldObjPt, $true, $false, $false, $false, $false, $true, 1, $false, $newObjPt, 2) | Out-Null

# 3) Objetivos (EN, italic) paragraph: drop the <w:t> entirely, keep the empty
#    italic run.
Set-ParagraphRuns "To introduce basic concepts on the influence" '<w:p><w:r><w:rPr><w:i/></w:rPr></w:r></w:p>' | Out-Null

# 4) Docente(s) list: append a second responsible professor, separated with a
#    line break, as its own run.
$docenteXml = '<w:p><w:r><w:t>5840897 - Clodoaldo Saron</w:t><w:br/></w:r><w:r><w:t>1033242 - F' + [char]0x00E1 + 'bio Herbst Florenzano</w:t></w:r></w:p>'
Set-ParagraphRuns "5840897 - Clodoaldo Saron" $docenteXml | Out-Null

# 5) Programa (PT) paragraph: replace with the updated syllabus text.
$oldProgPt = "Introdução: nomenclatura, arquitetura molecular e estrutura configuracional. Estado sólido: amorfo, cristalino e elastomérico. Estrutura e propriedades. Thermoplásticos: estrutura, propriedades e aplicações.  Elastômeros: estrutura, propriedades e aplicações. Resinas termorrígidas: estrutura, propriedades e aplicações.  Aditivos para polímeros: classes e aplicações. Avaliação"
$newProgPt = "Introdução: Desenvolvimento dos materiais poliméricos, organização da cadeia produtiva, formas de obtenção, nomenclatura, arquitetura molecular e estrutura configuracional. Estado sólido: amorfo, cristalino e elastomérico. Estrutura e propriedades. Thermoplásticos: estrutura, propriedades e aplicações. Elastômeros: estrutura, propriedades e aplicações. Resinas termorrígidas: estrutura, propriedades e aplicações. Propriedades mecânicas dos polímeros: comportamento à tração, impacto, flexão e fluência. Aditivos para polímeros: classes e aplicações. Viagem Didática complementar"
$d.Content.Find.Execute($oldProgPt, $true, $false, $false, $false, $false, $true, 1, $false, $newProgPt, 2) | Out-Null

# 6) Programa (EN, italic) paragraph: remove entirely (whole paragraph incl.
#    its paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Introduction: nomenclature, molecular architecture*") {
        $p.Range.Delete()
        break
    }
}

# 7) Bibliografia paragraph: replace with the updated reference list.
$oldBib = '1. L. C. SAWYER & D. T. GRUBB. Polymer Microscopy. New York: Chapman and Hall, 1987.2. H. G. ELLIAS Macromolecules -1-Structure and Properties, 2nd ed. Vol.1, New York, Plenum Press, 1984.3. J. M. G. COWIE Polymers: Chemistry and Physics of Modern Materials", New York: Chapman & Hall, 1998.4. C. A. HARPER Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill Inc, 1992.5. S. V. CANEVAROLO Jr. Ciência dos Polímeros. Editora Artiliber, 2001.6. E. B. MANO & L.C. MENDES Introdução a Polímeros. São Paulo: Ed. Edgard Blücher, 1999.7. E. B. MANO Polímeros como Materiais de Engenharia. São Paulo: Ed. Edgard Blücher, 1991.8. E. B. MANO & L. C. MENDES Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000.'
$newBib = 'SIMAL, A. L. Estrutura e Propriedades dos Polímeros, EduFSCar, São Carlos, 2002.SPERLING, L. H. Introduction to Physical Polymer Science, New York, John Wiley & Sans, 1985.BRYDSON, J. A. Rubbery Materials and Their Compounds, Elsevier, London, 1988.Rabello, M. S. Aditivação de Polímeros, Artiliber, São Paulo, 2004.HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill Inc, 1992 S. V. CANEVAROLO Jr. Técnicas de Caracterização de Polímeros. São Paulo: Editora Artliber, 2005. MANRICH, S. Processamento de Termoplásticos. Editora Artliber, 2005. NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997. MANO, E. B.; MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000. TURI, E. A. Thermal Characterization of Polymeric Materials. New York: Academic Press, 1981. NAVARRO, R.F. Fundamentos de Reologia de Polímeros. Editora da Universidade de Caxias do Sul, 1997.MANO, E. B.; MENDES, L. C. Identificação de Plásticos, Borrachas e Fibras. Ed. Edgard Blücher, 2000'
$d.Content.Find.Execute($oldBib, $true, $false, $false, $false, $false, $true, 1, $false, $newBib, 2) | Out-Null

Write-Host "Edits applied."
